# Updates the "cryptos" list (Sheet1) with refreshed prices / 1h volume
# percentages, matching the Sun Aug 18 19:24:15 UTC 2024 GitHub Actions
# data refresh. Row 16 and 17 also swap places (WrappedEther <-> ShibaInu)
# because the ranking order changed.
#
# Note: several "Price" (column D) values are plain numbers written as
# text in the workbook (e.g. "0.999", "6.71"). Assigning such a string
# straight to .Value lets Excel's type-inference turn it into a real
# number, which would not match the source data (still text). Prefixing
# with a single quote forces Excel to keep it as text (like typing '0.999
# into a cell); resetting the cell style back to Normal afterwards clears
# the "quote prefix" formatting flag that operation leaves behind, so the
# cell ends up identical to a plain text cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.628.61'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '2.649.93'
$ws.Range("E3").Value = '  +1.71%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''537.35'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("D6").Value = '''146.13'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +3.57%  '
$ws.Range("D7").Value = '''0.998'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '''0.573'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  +0.80%  '
$ws.Range("D9").Value = '''6.71'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +4.04%  '
$ws.Range("D10").Value = '''0.103'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("E11").Value = '  +1.32%  '
$ws.Range("E12").Value = '  -0.29%  '
$ws.Range("D13").Value = '3.112.42'
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").Value = '59.461.88'
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("D15").Value = '''21.31'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +3.62%  '

# Row 16 / 17 swap ranking places: WrappedEther <-> ShibaInu
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '''0.0000135'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +1.07%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.618.14'
$ws.Range("E17").Value = '  +0.30%  '

$ws.Range("D18").Value = '''340.23'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -0.85%  '
$ws.Range("D19").Value = '''4.41'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +1.38%  '
$ws.Range("D20").Value = '''10.33'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +2.08%  '
$ws.Range("E21").Value = '  -2.44%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").Value = '''66.73'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -1.13%  '
$ws.Range("E24").Value = '  +2.15%  '
$ws.Range("D25").Value = '''0.165'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  -1.07%  '
$ws.Range("E26").Value = '  -0.23%  '
$ws.Range("D27").Value = '''7.30'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +1.45%  '
$ws.Range("D28").Value = '0.0₃0745'
$ws.Range("E28").Value = '  +0.83%  '
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").Value = '''18.93'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +0.80%  '
$ws.Range("D33").Value = '''151.39'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  +1.35%  '
$ws.Range("D34").Value = '''4.01'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +0.74%  '
$ws.Range("D35").Value = '''1.14'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +1.98%  '
$ws.Range("D36").Value = '''0.847'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +2.69%  '
$ws.Range("D37").Value = '''0.837'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("D38").Value = '''1.45'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -1.05%  '
$ws.Range("E39").Value = '  +1.55%  '
$ws.Range("D40").Value = '''286.59'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +4.58%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").Value = '''0.608'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +2.06%  '
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("E44").Value = '  +3.12%  '
$ws.Range("D45").Value = '''19.31'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +4.08%  '
$ws.Range("E46").Value = '  -1.22%  '
$ws.Range("D47").Value = '''0.0227'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +1.77%  '
$ws.Range("D48").Value = '1.968.32'
$ws.Range("E48").Value = '  +1.08%  '
$ws.Range("D49").Value = '''4.57'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +1.28%  '
$ws.Range("D50").Value = '''18.31'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("D51").Value = '''111.12'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +0.14%  '
